# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the latest generated counts.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 143
    4  = 1399
    9  = 134
    11 = 346
    12 = 343
    13 = 1831
    17 = 719
    19 = 352
    20 = 4377
    22 = 320
    23 = 1183
    26 = 738
    28 = 387
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
